# Removendo o campo data do acervo arte gráfica
# The "Data" column (column G) is deleted entirely, shifting all
# subsequent columns one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G:G").Delete() | Out-Null

# Reflect the resulting selection left behind in the saved file.
$ws.Range("I10").Select() | Out-Null
